$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) writes to remain plain text even when the value
# looks numeric (e.g. "0.999", "1.00"), matching the inlineStr cells in the
# source file. NumberFormat "@" prevents Excel's auto-number coercion; the
# style is reset back to Normal afterwards so no stray style index sticks to
# the cell (keeps cell formatting identical to before the edit).
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "51.441.69"
$ws.Range("E2").Value = "  -0.41%  "

Set-TextValue $ws.Range("D3") "2.924.02"
$ws.Range("E3").Value = "  +0.87%  "

$ws.Range("E4").Value = "  -0.15%  "

Set-TextValue $ws.Range("D5") "361.63"
$ws.Range("E5").Value = "  +1.93%  "

Set-TextValue $ws.Range("D6") "104.15"
$ws.Range("E6").Value = "  -4.30%  "

$ws.Range("E7").Value = "  -2.86%  "

Set-TextValue $ws.Range("D8") "0.999"
$ws.Range("E8").Value = "  -0.22%  "

Set-TextValue $ws.Range("D9") "0.593"
$ws.Range("E9").Value = "  -4.60%  "

Set-TextValue $ws.Range("D10") "37.15"
$ws.Range("E10").Value = "  -4.13%  "

$ws.Range("E11").Value = "  +1.70%  "

Set-TextValue $ws.Range("D12") "0.0840"
$ws.Range("E12").Value = "  -3.00%  "

Set-TextValue $ws.Range("D13") "18.65"
$ws.Range("E13").Value = "  -3.79%  "

Set-TextValue $ws.Range("D14") "3.385.32"

Set-TextValue $ws.Range("D15") "7.39"
$ws.Range("E15").Value = "  -4.28%  "

Set-TextValue $ws.Range("D16") "2.921.90"
$ws.Range("E16").Value = "  +0.11%  "

Set-TextValue $ws.Range("D17") "0.968"
$ws.Range("E17").Value = "  -0.39%  "

Set-TextValue $ws.Range("D18") "51.344.83"
$ws.Range("E18").Value = "  -0.64%  "

$ws.Range("E19").Value = "  -1.53%  "

Set-TextValue $ws.Range("D20") "7.28"
$ws.Range("E20").Value = "  -2.97%  "

Set-TextValue $ws.Range("D21") "13.13"
$ws.Range("E21").Value = "  -5.09%  "

Set-TextValue $ws.Range("D22") "0.0₃0949"
$ws.Range("E22").Value = "  -2.42%  "

Set-TextValue $ws.Range("D23") "68.65"
$ws.Range("E23").Value = "  -2.34%  "

Set-TextValue $ws.Range("D24") "260.93"
$ws.Range("E24").Value = "  -2.40%  "

Set-TextValue $ws.Range("D25") "2.69"
$ws.Range("E25").Value = "  -3.51%  "

$ws.Range("E26").Value = "  -4.57%  "

$ws.Range("B27").Value = "Dai"
$ws.Range("C27").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue $ws.Range("D27") "1.00"
$ws.Range("E27").Value = "  +0.08%  "

$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue $ws.Range("D28") "26.22"
$ws.Range("E28").Value = "  -1.94%  "

$ws.Range("B29").Value = "Filecoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws.Range("D29") "7.35"
$ws.Range("E29").Value = "  -3.38%  "

$ws.Range("B30").Value = "Hedera"
$ws.Range("C30").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws.Range("D30") "0.110"
$ws.Range("E30").Value = "  +4.46%  "

$ws.Range("B31").Value = "RenderToken"
$ws.Range("C31").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D31") "6.22"
$ws.Range("E31").Value = "  +3.35%  "

$ws.Range("B32").Value = "Cosmos"
$ws.Range("C32").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue $ws.Range("D32") "10.04"
$ws.Range("E32").Value = "  -3.84%  "

$ws.Range("E33").Value = "  -2.21%  "

Set-TextValue $ws.Range("D34") "35.13"
$ws.Range("E34").Value = "  -5.18%  "

Set-TextValue $ws.Range("D35") "51.24"
$ws.Range("E35").Value = "  -1.74%  "

$ws.Range("E36").Value = "  +0.23%  "

Set-TextValue $ws.Range("D37") "0.0426"
$ws.Range("E37").Value = "  -3.18%  "

Set-TextValue $ws.Range("D38") "2.83"
$ws.Range("E38").Value = "  +4.84%  "

$ws.Range("E39").Value = "  +0.11%  "

Set-TextValue $ws.Range("D40") "17.12"
$ws.Range("E40").Value = "  -5.56%  "

Set-TextValue $ws.Range("D41") "1.88"
$ws.Range("E41").Value = "  -5.33%  "

$ws.Range("E42").Value = "  -3.73%  "

Set-TextValue $ws.Range("D43") "22.84"
$ws.Range("E43").Value = "  +0.76%  "

Set-TextValue $ws.Range("D44") "119.85"
$ws.Range("E44").Value = "  +0.90%  "

$ws.Range("E45").Value = "  -1.32%  "

Set-TextValue $ws.Range("D46") "2.087.96"
$ws.Range("E46").Value = "  -1.53%  "

Set-TextValue $ws.Range("D47") "3.22"
$ws.Range("E47").Value = "  -5.74%  "

$ws.Range("E48").Value = "  -7.27%  "

Set-TextValue $ws.Range("D49") "3.220.97"
$ws.Range("E49").Value = "  +0.67%  "

Set-TextValue $ws.Range("D50") "0.240"
$ws.Range("E50").Value = "  -3.88%  "

Set-TextValue $ws.Range("D51") "0.0316"
$ws.Range("E51").Value = "  -7.29%  "
